$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Rewrite the three existing "DYOR" list items that survive (SOL -> LUNA
#    text, TRX -> SOL text, XRP -> TRX text). Each paragraph keeps its
#    existing Paragraphedeliste / numPr formatting; we only touch the text.
# ---------------------------------------------------------------------------

function Set-ParaText($para, [string]$text) {
    $r = $para.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

$pSol = $d.Paragraphs.Item(8)
Set-ParaText $pSol "LUNA: like an ecosystem, Terra & Luna. Terra is the stablecoin. The more people will buy Terre the more Luna will increase in price because for minting Terra it needs to burn Luna."

$pTrx = $d.Paragraphs.Item(9)
Set-ParaText $pTrx "SOL: scaling, Pos & PoH, High TPS, low fees."

$pXrp = $d.Paragraphs.Item(10)
Set-ParaText $pXrp "TRX: Ethereum-like, DPOS instead: High TPS and low fees."

# ---------------------------------------------------------------------------
# 2) Insert four brand-new list items (XRP, ADA, "Know the 20-1st...", WIN)
#    right after the paragraph that used to hold "XRP" (now "TRX: ..."),
#    and before the "Miscellaneous" heading. Inserting after a paragraph
#    that already carries the Paragraphedeliste / numPr(ilvl=0,numId=5)
#    formatting makes the new paragraphs inherit the same list formatting.
# ---------------------------------------------------------------------------

$pXrp.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs.Item(11)
Set-ParaText $pNew1 "XRP: Unique Node List protocol, not very decentralized, stated as a security rather than a currency."

$pNew1.Range.InsertParagraphAfter()
$pNew2 = $d.Paragraphs.Item(12)
Set-ParaText $pNew2 "ADA: is like ETH 2.0 but deflationary: 21 billion tokens. Process of verification for different area: education, retail, agriculture, medical, finance, government."

$pNew2.Range.InsertParagraphAfter()
$pNew3 = $d.Paragraphs.Item(13)
Set-ParaText $pNew3 "Know the 20-1st crypto by market cap."

$pNew3.Range.InsertParagraphAfter()
$pNew4 = $d.Paragraphs.Item(14)
# Include the trailing space up front: adding the bookmark exactly at the
# paragraph's text end (i.e. with nothing after it) mis-places the
# bookmark, so the space must already be present before we add it.
Set-ParaText $pNew4 "WIN: Oracle for TRX network "

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark off the "Miscellaneous" heading and onto
#    the new "WIN: ..." paragraph (right after the text, matching where
#    Word leaves the caret after the last edit).
# ---------------------------------------------------------------------------

$pMisc = $d.Paragraphs.Item(15)
$oldBookmark = $pMisc.Range.Bookmarks("_GoBack")
$oldBookmark.Delete()

$winTextLen = ("WIN: Oracle for TRX network").Length
$bmPos = $pNew4.Range.Start + $winTextLen
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
